$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("runsettings TOT")
$ws3 = $wb.Worksheets.Item("spcode TOT")

# --- Sheet2 "runsettings TOT": add rows 3,4,5 duplicating row 2 with different runno / runname ---
$runnames = @("ufs2", "ufs3", "ufs4")
for ($i = 0; $i -lt 3; $i++) {
    $r = 3 + $i
    $ws2.Cells.Item($r, 1).Value = $i + 2
    $ws2.Cells.Item($r, 2).Value = $runnames[$i]
    $ws2.Cells.Item($r, 3).Value = "policyData1250f"
    $ws2.Cells.Item($r, 4).Value = "alloc202012"
    $ws2.Cells.Item($r, 5).Value = 202012
    $ws2.Cells.Item($r, 6).Value = 12
    $ws2.Cells.Item($r, 7).Value = "qxtables.xlsx"
    $ws2.Cells.Item($r, 8).Value = "wtables.xlsx"
    $ws2.Cells.Item($r, 9).Value = "vtables.xlsx"
    $ws2.Cells.Item($r, 10).Value = 1
    $ws2.Cells.Item($r, 11).Value = 0.0
    $ws2.Cells.Item($r, 12).Value = 0.0
    $ws2.Cells.Item($r, 13).Value = 1
    $ws2.Cells.Item($r, 14).Value = 1
    $ws2.Cells.Item($r, 15).Value = 1
    $ws2.Cells.Item($r, 16).Value = "exptables.xlsx"
    $ws2.Cells.Item($r, 17).Value = "chtables.xlsx"
    $ws2.Cells.Item($r, 18).Value = "fundtables.xlsx"
    $ws2.Cells.Item($r, 19).Value = "ul1"
}

# --- Sheet3 "spcode TOT": add rows 7-21, three blocks of 5 rows each (spc=2,3,4) ---
$spcodes = @("ul1", "SA85", "SA85")
$wxbasisVals = @(
    @("prot", "prot", "prot", "ann", "prot"),
    @("ASSA", "ASSA", "ASSA", "ASSA", "ASSA"),
    @("prot", "prot", "prot", "ann", "prot")
)
$vtbasisVals = @("v5", "v5", "v2")
$aVals = @(10, 11, 20, 30, 40)

$row = 7
for ($b = 0; $b -lt 3; $b++) {
    $spc = $b + 2
    for ($k = 0; $k -lt 5; $k++) {
        $ws3.Cells.Item($row, 1).Value = $aVals[$k]
        $ws3.Cells.Item($row, 2).Value = $spc
        $ws3.Cells.Item($row, 3).Value = $spcodes[$b]
        $ws3.Cells.Item($row, 4).Value = $wxbasisVals[$b][$k]
        $ws3.Cells.Item($row, 5).Value = $vtbasisVals[$b]
        $ws3.Cells.Item($row, 6).Value = "ufs1"
        $ws3.Cells.Item($row, 7).Value = $spcodes[$b]
        $ws3.Cells.Item($row, 8).Value = $wxbasisVals[$b][$k]
        $ws3.Cells.Item($row, 9).Value = $vtbasisVals[$b]
        $ws3.Cells.Item($row, 10).Value = "ufs1"
        $ws3.Cells.Item($row, 11).Value = "zero"
        $row++
    }
}
